$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.032364726066589
$ws.Range("B1").Value = 2.00324273109436
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.793906211853027
$ws.Range("E1").Value = 1.160421252250671
